$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 55571
$ws.Range("I11").Value = 55571
$ws.Range("K11").Value = 55571
$ws.Range("M11").Value = -55431

$ws.Range("H19").Value = 1563.5
$ws.Range("I19").Value = 1357.8334
$ws.Range("K19").Value = 1357.8334
$ws.Range("M19").Value = -1182.8334

$ws.Range("H21").Value = 9480.083000000001
$ws.Range("I21").Value = 9480.083000000001
$ws.Range("K21").Value = 9480.083000000001
$ws.Range("M21").Value = -9012.083000000001

$ws.Range("H23").Value = 9480.083000000001
$ws.Range("I23").Value = 9480.083000000001
$ws.Range("K23").Value = 9480.083000000001
$ws.Range("M23").Value = -9246.083000000001

$ws.Range("H28").Value = 1495.3334
$ws.Range("J28").Value = 4006
$ws.Range("L28").Value = 4006
$ws.Range("N28").Value = -4976

$ws.Range("H47").Value = 1234
$ws.Range("I47").Value = 1234
$ws.Range("K47").Value = 1234
$ws.Range("M47").Value = -262

$ws.Range("H54").Value = 5555
$ws.Range("I54").Value = 5555
$ws.Range("K54").Value = 5555
$ws.Range("M54").Value = -5069

$ws.Range("H80").Value = 4308.5454
$ws.Range("I80").Value = 794.5
$ws.Range("J80").Value = 5089.4443
$ws.Range("K80").Value = 2383.5
$ws.Range("L80").Value = 15268.3329
$ws.Range("M80").Value = -1385.5
$ws.Range("N80").Value = -17264.3329

$ws.Range("H83").Value = 4308.5454
$ws.Range("I83").Value = 794.5
$ws.Range("J83").Value = 5089.4443
$ws.Range("K83").Value = 7150.5
$ws.Range("L83").Value = 45804.9987
$ws.Range("M83").Value = -2158.5
$ws.Range("N83").Value = -55788.9987

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 821.9818
$ws.Range("I32").Value = 821.9818
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 821.9818
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -534.9818
$ws.Range("N32").ClearContents()

$ws.Range("H74").Value = 7719915.5
$ws.Range("I74").Value = 3708343.5
$ws.Range("J74").Value = 27777776
$ws.Range("K74").Value = 3708343.5
$ws.Range("L74").Value = 27777776
$ws.Range("M74").Value = -3707469.5
$ws.Range("N74").Value = -27779524

$ws.Range("H77").Value = 7719915.5
$ws.Range("I77").Value = 3708343.5
$ws.Range("J77").Value = 27777776
$ws.Range("K77").Value = 18541717.5
$ws.Range("L77").Value = 138888880
$ws.Range("M77").Value = -18537349.5
$ws.Range("N77").Value = -138897616

$ws.Range("H110").Value = 2993
$ws.Range("I110").Value = 670.2857
$ws.Range("K110").Value = 670.2857
$ws.Range("M110").Value = 1374.7143

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 16691
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3238.8333
$ws.Range("I107").Value = 1897
$ws.Range("K107").Value = 1897
$ws.Range("M107").Value = 23

$ws.Range("H134").Value = 3227588.2
$ws.Range("I134").Value = 1490.8276
$ws.Range("J134").Value = 50006000
$ws.Range("K134").Value = 4472.4828
$ws.Range("L134").Value = 150018000
$ws.Range("M134").Value = -1937.4828
$ws.Range("N134").Value = -150023070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 38999
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 38999
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 38999
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -40017

$ws.Range("H132").Value = 1797.2
$ws.Range("I132").Value = 1797.2
$ws.Range("K132").Value = 5391.6
$ws.Range("M132").Value = -2861.6

$ws.Range("H136").Value = 33757.832
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 33757.832
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 101273.496
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -106373.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1624.25
$ws.Range("I7").Value = 1624.25
$ws.Range("K7").Value = 1624.25
$ws.Range("M7").Value = -1512.25

$ws.Range("H9").Value = 5453.857
$ws.Range("J9").Value = 8742
$ws.Range("L9").Value = 8742
$ws.Range("N9").Value = -9190

$ws.Range("H10").Value = 2528.5715
$ws.Range("I10").Value = 2066.6667
$ws.Range("J10").Value = 2875
$ws.Range("K10").Value = 2066.6667
$ws.Range("L10").Value = 2875
$ws.Range("M10").Value = -1926.6667
$ws.Range("N10").Value = -3155

$ws.Range("H46").Value = 1794.6097
$ws.Range("I46").Value = 1424.3334
$ws.Range("K46").Value = 1424.3334
$ws.Range("M46").Value = -1236.3334

$ws.Range("H126").Value = 1624.25
$ws.Range("I126").Value = 1624.25
$ws.Range("K126").Value = 4872.75
$ws.Range("M126").Value = -2402.75

$ws.Range("H136").Value = 250003500
$ws.Range("I136").Value = 3499.5
$ws.Range("J136").Value = 500003500
$ws.Range("K136").Value = 10498.5
$ws.Range("L136").Value = 1500010500
$ws.Range("M136").Value = -7948.5
$ws.Range("N136").Value = -1500015600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2322.6667
$ws.Range("I126").Value = 1666.8
$ws.Range("J126").Value = 2791.1428
$ws.Range("K126").Value = 5000.4
$ws.Range("L126").Value = 8373.428400000001
$ws.Range("M126").Value = -2530.4
$ws.Range("N126").Value = -13313.4284

$ws.Range("H136").Value = 1924.8334
$ws.Range("I136").Value = 1828.5714
$ws.Range("J136").Value = 2059.6
$ws.Range("K136").Value = 5485.7142
$ws.Range("L136").Value = 6178.799999999999
$ws.Range("M136").Value = -2935.7142
$ws.Range("N136").Value = -11278.8

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
